# Apply updated dSF ("F") column values.
# This reflects a repull/recalculation of the underlying data (per commit
# message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 4
    "F4"  = 1
    "F8"  = 1
    "F12" = -10
    "F18" = -6
    "F19" = -3
    "F24" = -1
    "F25" = -4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
